$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 31 de Mayo de 2020 a las 06:05"

# Pakistan row (row 21) - updated stats
$ws.Range("B21").Value = 69496
$ws.Range("C21").Value = 3039
$ws.Range("D21").Value = 25271
$ws.Range("E21").Value = 42742
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 1483

# Rows 53-55: Kazajistan moves above Barein/Oman (sorted by total cases)
$ws.Range("A53").Value = "Kazajistan"
$ws.Range("B53").Value = 10858
$ws.Range("C53").Value = 476
$ws.Range("D53").Value = 5220
$ws.Range("E53").Value = 5600
$ws.Range("H53").Value = 38

$ws.Range("A54").Value = "Barein"
$ws.Range("B54").Value = 10793
$ws.Range("C54").Value = 0
$ws.Range("D54").Value = 5826
$ws.Range("E54").Value = 4950
$ws.Range("H54").Value = 17

$ws.Range("A55").Value = "Oman"
$ws.Range("B55").Value = 10423
$ws.Range("C55").Value = 0
$ws.Range("D55").Value = 2396
$ws.Range("E55").Value = 7985
$ws.Range("H55").Value = 42

# Rows 200-201: Belice moves above Santa Lucia
$ws.Range("A200").Value = "Belice"
$ws.Range("D200").Value = 16
$ws.Range("H200").Value = 2

$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("D201").Value = 18
$ws.Range("H201").Value = 0

# Rows 213-214: Islas Virgenes Britanicas moves above Papua Nueva Guinea
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
